$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns G and H
$ws.Range("G1").Value = "ci.lower"
$ws.Range("H1").Value = "ci.upper"

$ws.Range("G2").Value = -0.0619076051069637
$ws.Range("H2").Value = 0.34115075752665
$ws.Range("G3").Value = -0.0523276801642743
$ws.Range("H3").Value = 0.0848699822433146
$ws.Range("G4").Value = -0.0731788694076288
$ws.Range("H4").Value = 0.0723196290497396
$ws.Range("G5").Value = -0.169389800520079
$ws.Range("H5").Value = 0.274732403905662
$ws.Range("G6").Value = -0.217410778806589
$ws.Range("H6").Value = 0.214858018468764
$ws.Range("G7").Value = -0.162127835312166
$ws.Range("H7").Value = 0.262954261700384
$ws.Range("G8").Value = -0.226731466986582
$ws.Range("H8").Value = 0.224069266430392
$ws.Range("G9").Value = -0.0617481755265165
$ws.Range("H9").Value = 0.0775897063870921
$ws.Range("G10").Value = 0.152802639089172
$ws.Range("H10").Value = 0.201314135751746
$ws.Range("G11").Value = -0.0074050499760908
$ws.Range("H11").Value = 0.0408065924130202
$ws.Range("G15").Value = -0.0688198396151368
$ws.Range("H15").Value = 0.0863798511914907
$ws.Range("G16").Value = -0.213225803022513
$ws.Range("H16").Value = 0.267632317050905
$ws.Range("G17").Value = -0.194258792173918
$ws.Range("H17").Value = 0.250759234373449
$ws.Range("G18").Value = -0.144639258099479
$ws.Range("H18").Value = 0.146357738815258
$ws.Range("G19").Value = -0.448138532860785
$ws.Range("H19").Value = 0.453462933973163
$ws.Range("G20").Value = -0.429716036937527
$ws.Range("H20").Value = 0.434821557613179
$ws.Range("G21").Value = -0.0074050499760908
$ws.Range("H21").Value = 0.0408065924130202
$ws.Range("G22").Value = -0.123496351053033
$ws.Range("H22").Value = 0.155179412774184
$ws.Range("G23").Value = -0.0229432055698443
$ws.Range("H23").Value = 0.126431832514251
$ws.Range("G24").Value = -0.382631066432205
$ws.Range("H24").Value = 0.480795292264233
$ws.Range("G25").Value = -0.0199386179492627
$ws.Range("H25").Value = 0.127833981672672
$ws.Range("G26").Value = -0.380564418284175
$ws.Range("H26").Value = 0.483354261331932

# Rows 12-14 have no CI values (NA in the source data); materialize
# them as empty cells (matching the existing empty C/D/E/F pattern
# in those rows) without introducing any new cell style.
$ws.Range("G12:H14").Style = "Normal"
